$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.76
$ws.Range("O2").Value = 1.36
$ws.Range("S2").Value = 3.55
$ws.Range("U2").Value = 2.12
$ws.Range("W2").Value = 1.56
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 60
$ws.Range("AB2").Value = 10.5
$ws.Range("AC2").Value = 8.6
$ws.Range("AD2").Value = 13.5
$ws.Range("AG2").Value = 12.5
$ws.Range("AJ2").Value = 48
$ws.Range("AO2").Value = 40
$ws.Range("N3").Value = 3.4
$ws.Range("G4").Value = 3.85
$ws.Range("I4").Value = 2.32
$ws.Range("J4").Value = 3.7
$ws.Range("L4").Value = 1.27
$ws.Range("M4").Value = 1.04
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 1.49
$ws.Range("S4").Value = 2.6
$ws.Range("T4").Value = 1.58
$ws.Range("U4").Value = 2.34
$ws.Range("V4").Value = 1.76
$ws.Range("W4").Value = 1.35
$ws.Range("Y4").Value = 18
$ws.Range("Z4").Value = 22
$ws.Range("AA4").Value = 38
$ws.Range("AB4").Value = 24
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 15.5
$ws.Range("AE4").Value = 30
$ws.Range("AF4").Value = 38
$ws.Range("AH4").Value = 22
$ws.Range("AI4").Value = 42
$ws.Range("K5").Value = 4.9
$ws.Range("L5").Value = 1.28
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 2.16
$ws.Range("O5").Value = 1.23
$ws.Range("R5").Value = 1.45
$ws.Range("S5").Value = 2.76
$ws.Range("T5").Value = 1.01
$ws.Range("U5").Value = 1.01
$ws.Range("V5").Value = 1.17
$ws.Range("W5").Value = 2.38
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000
$ws.Range("F6").Value = 1.33
$ws.Range("G6").Value = 1.42
$ws.Range("I6").Value = 12
$ws.Range("J6").Value = 5.3
$ws.Range("K6").Value = 6.4
$ws.Range("P6").Value = 2.52
$ws.Range("Q6").Value = 1.51
$ws.Range("H7").Value = 5.3
$ws.Range("Q7").Value = 1.64
$ws.Range("R7").Value = 1.58
$ws.Range("Z7").Value = 55
$ws.Range("AE7").Value = 70
$ws.Range("AL7").Value = 34
$ws.Range("AO7").Value = 70
$ws.Range("Q8").Value = 2.02
$ws.Range("N9").Value = 5.6
$ws.Range("Q9").Value = 1.62
$ws.Range("R9").Value = 1.61
$ws.Range("U9").Value = 2.66
$ws.Range("N10").Value = 3.7
$ws.Range("F11").Value = 1.46
$ws.Range("I11").Value = 8.199999999999999
$ws.Range("N11").Value = 5.8
$ws.Range("X11").Value = 34
$ws.Range("Y11").Value = 38
$ws.Range("Z11").Value = 75
$ws.Range("AB11").Value = 14
$ws.Range("AD11").Value = 36
$ws.Range("AE11").Value = 120
$ws.Range("AH11").Value = 23
$ws.Range("AI11").Value = 95
$ws.Range("AL11").Value = 32
$ws.Range("AO11").Value = 100
$ws.Range("F12").Value = 1.96
$ws.Range("G12").Value = 1.98
$ws.Range("H12").Value = 4.5
$ws.Range("I12").Value = 4.8
$ws.Range("P12").Value = 1.79
$ws.Range("Q12").Value = 2.22
$ws.Range("Z12").Value = 34
$ws.Range("AN12").Value = 17
